# Applies updated '想去人数' (F column) and '最低票价' (G column) stats
# to the 展览, 演出, and 全部类型 sheets, matching the refreshed data
# generated for the gh-pages output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 2779
$ws.Range("G6").Value = 75
$ws.Range("F7").Value = 1655
$ws.Range("G7").Value = 70
$ws.Range("F8").Value = 7498
$ws.Range("F10").Value = 7689
$ws.Range("F12").Value = 38
$ws.Range("F13").Value = 17
$ws.Range("F14").Value = 6272
$ws.Range("F15").Value = 3286
$ws.Range("F17").Value = 3657
$ws.Range("F20").Value = 22
$ws.Range("F21").Value = 38
$ws.Range("F22").Value = 6
$ws.Range("F25").Value = 290
$ws.Range("F26").Value = 297
$ws.Range("F27").Value = 3682
$ws.Range("F30").Value = 931
$ws.Range("F32").Value = 1325
$ws.Range("F35").Value = 2642
$ws.Range("F36").Value = 1573
$ws.Range("F38").Value = 29
$ws.Range("F39").Value = 34
$ws.Range("F40").Value = 3355
$ws.Range("F41").Value = 198
$ws.Range("F44").Value = 903
$ws.Range("F45").Value = 492
$ws.Range("F46").Value = 1312
$ws.Range("F47").Value = 231
$ws.Range("F48").Value = 531
$ws.Range("F49").Value = 602

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G4").Value = 108
$ws.Range("F5").Value = 236
$ws.Range("F6").Value = 45
$ws.Range("F7").Value = 20
$ws.Range("F8").Value = 40
$ws.Range("F10").Value = 32
$ws.Range("F12").Value = 100
$ws.Range("F13").Value = 21
$ws.Range("F18").Value = 12

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("G5").Value = 108
$ws.Range("F8").Value = 2779
$ws.Range("G8").Value = 75
$ws.Range("F9").Value = 1655
$ws.Range("G9").Value = 70
$ws.Range("F10").Value = 236
$ws.Range("F11").Value = 45
$ws.Range("F12").Value = 40
$ws.Range("F13").Value = 7498
$ws.Range("F14").Value = 7689
$ws.Range("F16").Value = 38
$ws.Range("F17").Value = 6272
$ws.Range("F18").Value = 3286
$ws.Range("F19").Value = 3657
$ws.Range("F21").Value = 38
$ws.Range("F23").Value = 32
$ws.Range("F24").Value = 290
$ws.Range("F25").Value = 297
$ws.Range("F26").Value = 3682
$ws.Range("F30").Value = 931
$ws.Range("F32").Value = 1325
$ws.Range("F35").Value = 2642
$ws.Range("F36").Value = 1573
$ws.Range("F38").Value = 29
$ws.Range("F39").Value = 34
$ws.Range("F40").Value = 3355
$ws.Range("F41").Value = 198
$ws.Range("F44").Value = 903
$ws.Range("F45").Value = 492
$ws.Range("F46").Value = 1312
$ws.Range("F47").Value = 231
$ws.Range("F48").Value = 531
$ws.Range("F49").Value = 602
